$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E header + existing rows (2,3) description ---
$ws.Range("E1").Value = "Description"
$ws.Range("E2").Value = "Valid credentials"
$ws.Range("E3").Value = "Invalid Username and Password"

# --- New row 4: valid username, invalid password ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Jonnny_Doe"
$ws.Range("C4").Value = "InvalidPassword"
$ws.Range("D4").Formula = "=FALSE"
$ws.Range("E4").Value = "Valid Username and Invalid Password"

# --- New row 5: invalid username, valid password ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "InvalidUserName2"
$ws.Range("C5").Value = "bkxEGLmduZeG4UaiSLQn"
$ws.Range("D5").Formula = "=FALSE"
$ws.Range("E5").Value = "Invalid Username and Valid Password"

# --- Row 2 changes: password becomes a new (slightly different) value, IsValid becomes FALSE ---
$ws.Range("C2").Value = "bkxEGLmduZeG4UaiSLQn1"
$ws.Range("D2").Formula = "=FALSE"

# --- Resize columns to fit their new content ---
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null

# --- Match the recorded selection ---
$ws.Range("D2").Select()
